$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.743.65"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "2.084.21"
$ws.Range("E3").Value = "  +0.26%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.90"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.87%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  +1.86%  "
$ws.Range("E10").Value = "  +0.53%  "
$ws.Range("E11").Value = "  +2.81%  "
$ws.Range("D12").Value = "2.390.52"
$ws.Range("E12").Value = "  +0.58%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.76"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.98%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.23"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.22%  "
$ws.Range("E15").Value = "  +1.76%  "
$ws.Range("E16").Value = "  +1.19%  "
$ws.Range("D17").Value = "2.104.46"
$ws.Range("E17").Value = "  +1.18%  "
$ws.Range("D18").Value = "37.726.87"
$ws.Range("E18").Value = "  +0.07%  "
$ws.Range("E19").Value = "  -0.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.79"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.61%  "
$ws.Range("D21").Value = "0.0₃0846"
$ws.Range("E21").Value = "  +3.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.36"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.23%  "
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("E24").Value = "  -0.49%  "
$ws.Range("E25").Value = "  +1.26%  "
$ws.Range("E26").Value = "  +7.82%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "171.16"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.65%  "
$ws.Range("E28").Value = "  -1.89%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.43"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.94%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.57"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.93%  "
$ws.Range("E31").Value = "  +2.20%  "
$ws.Range("E32").Value = "  +2.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0634"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.49%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.69"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.80%  "
$ws.Range("E35").Value = "  -0.49%  "
$ws.Range("E36").Value = "  +0.98%  "
$ws.Range("E37").Value = "  -0.15%  "
$ws.Range("E38").Value = "  +0.23%  "
$ws.Range("E39").Value = "  +0.30%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0981"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.37%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.35"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +11.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "99.08"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.34%  "
$ws.Range("E43").Value = "  +2.56%  "
$ws.Range("D45").Value = "1.452.99"
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("E46").Value = "  -0.42%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.16"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.67%  "
$ws.Range("E48").Value = "  +1.23%  "
$ws.Range("E49").Value = "  -0.87%  "
$ws.Range("D51").Value = "2.276.47"
$ws.Range("E51").Value = "  +0.36%  "
